$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D-column cells whose new values would otherwise
# be auto-parsed by Excel as numbers (losing the source's text formatting),
# so they round-trip as plain text exactly like the scraped values.
$textForceCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values scraped by the GitHub Actions job.
$ws.Range("D2").Value = "30.497.25"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "1.937.21"
$ws.Range("E3").Value = "  +4.46%  "
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "240.70"
$ws.Range("E5").Value = "  +3.00%  "
$ws.Range("D6").Value = "0.9993"
$ws.Range("D7").Value = "0.4767"
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").Value = "44.68"
$ws.Range("E8").Value = "  +3.15%  "
$ws.Range("D9").Value = "0.2882"
$ws.Range("E9").Value = "  +4.46%  "
$ws.Range("D10").Value = "0.06648"
$ws.Range("E10").Value = "  +4.59%  "
$ws.Range("D11").Value = "107.81"
$ws.Range("E11").Value = "  +26.96%  "
$ws.Range("D12").Value = "19.06"
$ws.Range("E12").Value = "  +5.99%  "
$ws.Range("D13").Value = "1.928.26"
$ws.Range("E13").Value = "  +4.06%  "
$ws.Range("D14").Value = "0.07608"
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("D15").Value = "5.178"
$ws.Range("E15").Value = "  +3.63%  "
$ws.Range("D16").Value = "0.6628"
$ws.Range("E16").Value = "  +5.91%  "
$ws.Range("D17").Value = "305.90"
$ws.Range("E17").Value = "  +20.83%  "
$ws.Range("D18").Value = "30.514.55"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").Value = "13.05"
$ws.Range("E19").Value = "  +2.44%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.000007605"
$ws.Range("E20").Value = "  +3.38%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "2.167.55"
$ws.Range("E22").Value = "  +2.98%  "
$ws.Range("D23").Value = "5.291"
$ws.Range("E23").Value = "  +7.18%  "
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "6.319"
$ws.Range("E25").Value = "  +6.67%  "
$ws.Range("D26").Value = "9.331"
$ws.Range("E26").Value = "  +3.37%  "
$ws.Range("D27").Value = "168.30"
$ws.Range("E27").Value = "  +2.48%  "
$ws.Range("D28").Value = "20.87"
$ws.Range("E28").Value = "  +15.48%  "
$ws.Range("D29").Value = "2.059"
$ws.Range("E29").Value = "  +9.39%  "
$ws.Range("D30").Value = "0.1112"
$ws.Range("E30").Value = "  +8.44%  "
$ws.Range("D31").Value = "1.360"
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("D32").Value = "4.108"
$ws.Range("E32").Value = "  +1.39%  "
$ws.Range("D33").Value = "3.959"
$ws.Range("E33").Value = "  +2.92%  "
$ws.Range("E34").Value = "  +3.73%  "
$ws.Range("D35").Value = "0.7456"
$ws.Range("E35").Value = "  +6.33%  "
$ws.Range("D36").Value = "1.161"
$ws.Range("E36").Value = "  +2.42%  "
$ws.Range("D37").Value = "2.747"
$ws.Range("E37").Value = "  +2.09%  "
$ws.Range("D38").Value = "0.01965"
$ws.Range("E38").Value = "  +3.58%  "
$ws.Range("D39").Value = "2.694"
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("D40").Value = "2.041"
$ws.Range("E40").Value = "  +2.42%  "
$ws.Range("D41").Value = "0.8826"
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("D43").Value = "70.58"
$ws.Range("E43").Value = "  +11.17%  "
$ws.Range("D44").Value = "5.815"
$ws.Range("E44").Value = "  +5.22%  "
$ws.Range("D45").Value = "0.9990"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "0.4203"
$ws.Range("E46").Value = "  +3.04%  "
$ws.Range("D47").Value = "7.275"
$ws.Range("E47").Value = "  +1.12%  "
$ws.Range("D48").Value = "9.293"
$ws.Range("E48").Value = "  +8.72%  "
$ws.Range("D49").Value = "0.1217"
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("D50").Value = "34.93"
$ws.Range("E50").Value = "  +2.43%  "
$ws.Range("E51").Value = "  +2.36%  "
